$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("ae36fa6e-f222-4071-aec0-1c31450a4470", "Create Country", "PASSED", "03_28_2024_19_53_58", "03_28_2024_19_54_06", "PT8.3587541S"),
    @("4b08957a-055f-48e6-8ec5-01493c4c5f9e", "Delete Country", "PASSED", "03_28_2024_19_54_11", "03_28_2024_19_54_19", "PT8.0514124S"),
    @("2958d93b-e068-4570-8ae0-9059f8fff0cb", "Create country with parameters", "PASSED", "03_28_2024_19_54_22", "03_28_2024_19_54_31", "PT8.6318558S"),
    @("adf99497-3ed2-4980-9f65-2c598b736c6b", "Delete country with parameters", "PASSED", "03_28_2024_19_54_34", "03_28_2024_19_54_42", "PT7.9847116S")
)

$startRow = 70
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Range("D$row:F$row").HorizontalAlignment = -4108
}
